$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Collapse the whole "Task 2" question block into a single paragraph
#    that just reads "Task 3" (the grading-question sub-items, the
#    worked answers and the old standalone "Task 3" paragraph are all
#    removed), and bring the _GoBack bookmark along with it.
# ---------------------------------------------------------------------

# Find the paragraph that holds the lone "Task 2" heading text.
$task2Range = $d.Content
$task2Range.Find.ClearFormatting()
$found = $task2Range.Find.Execute("Task 2", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'Task 2' paragraph"
}
$task2ParaIndex = $task2Range.Paragraphs.Item(1).Index

# Find the paragraph that holds the old standalone "Task 3" heading text.
$task3Range = $d.Content
$task3Range.Find.ClearFormatting()
$found = $task3Range.Find.Execute("Task 3", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'Task 3' paragraph"
}
$task3ParaIndex = $task3Range.Paragraphs.Item(1).Index

# Delete every paragraph strictly between "Task 2" and "Task 3" plus the
# old "Task 3" paragraph itself, leaving the "Task 2" paragraph as the
# sole survivor (it becomes the new "Task 3" paragraph below). Paragraphs
# are always removed at the same index because later ones shift down.
$deleteCount = $task3ParaIndex - $task2ParaIndex
for ($i = 1; $i -le $deleteCount; $i++) {
    $d.Paragraphs.Item($task2ParaIndex + 1).Range.Delete()
}

# Rename the surviving paragraph's run text from "Task 2" to "Task 3".
# Scoping the Find/Replace to just that paragraph preserves the run's
# formatting (rFonts/sz/szCs) instead of recreating a plain run.
$taskPara = $d.Paragraphs.Item($task2ParaIndex)
$taskPara.Range.Find.ClearFormatting()
$taskPara.Range.Find.Execute("Task 2", $true, $true, $false, $false, $false, $true, 1, $false, "Task 3", 2) | Out-Null

# Word keeps a single "_GoBack" bookmark marking the last edit location;
# re-adding it here moves it from its old spot to the start of this
# paragraph, matching where Word would leave it after this edit.
$d.Bookmarks.Add("_GoBack", $d.Range($taskPara.Range.Start, $taskPara.Range.Start)) | Out-Null

# ---------------------------------------------------------------------
# 2. Re-save the "In your own words..." paragraph text in place so Word
#    consolidates its three split runs into one (no wording changes).
# ---------------------------------------------------------------------
$bugsText = "In your own words, explain how you implemented each task. Did you encounter any bugs? If so, how did you fix them? If you failed to complete any tasks, list them here and briefly explain why."
$bugsRange = $d.Content
$bugsRange.Find.ClearFormatting()
$found = $bugsRange.Find.Execute($bugsText, $true, $true, $false, $false, $false, $true, 1, $false, $bugsText, 2)
if (-not $found) {
    throw "Could not locate the 'In your own words...' paragraph"
}

# ---------------------------------------------------------------------
# 3. Re-save the "I used a struct..." sentence in place so Word merges
#    its many split/proofed runs into one (no wording changes); the
#    "For task 2, " lead-in keeps its own separate runs/formatting.
# ---------------------------------------------------------------------
$structText = "I used a struct and typecasting to pass the data to the function. I then used the random number generator to randomly pick with thread to run next. If thread 1 randomly chose a number that was not 1, the problem would call currentThread->Yield(). I then called a while loop to randomly shout output until the limit was reached. I learned how to used the random number generator in Nachos"
$structRange = $d.Content
$structRange.Find.ClearFormatting()
$found = $structRange.Find.Execute($structText, $true, $true, $false, $false, $false, $true, 1, $false, $structText, 2)
if (-not $found) {
    throw "Could not locate the 'I used a struct...' sentence"
}
